$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates: new reporting week (1/20/2025 - 1/26/2025) and the
# bulletin's volume/number counter (Number 3 -> Number 4).
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  4"
$ws.Range("C9").Value = "Report Covering the Week  1/20/2025  Through  1/26/2025"

# ---------------------------------------------------------------------------
# Helpers for writing into the crime-stats grid (rows 16-33, cols C:N).
# Numeric cells use style 14 (counts) / 15 (percent change); when a figure is
# unavailable the sheet instead shows literal text "0" or "***.*" drawn with
# style 13 (General format). Stable reference cells that keep these exact
# styles throughout this edit are used as format-paint sources so every
# touched cell ends up with the right style id (not a freshly minted one).
# ---------------------------------------------------------------------------
$numStyleSrc = $ws.Range("I16")   # style 14 - plain integer count
$pctStyleSrc = $ws.Range("L28")   # style 15 - percent-change number
$zeroTextSrc = $ws.Range("C20")   # style 13 - text "0"
$naTextSrc   = $ws.Range("E22")   # style 13 - text "***.*"

function Set-NumCell($ref, $value) {
    $ws.Range($ref).Value = $value
    $numStyleSrc.Copy() | Out-Null
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null
}

function Set-PctCell($ref, $value) {
    $ws.Range($ref).Value = $value
    $pctStyleSrc.Copy() | Out-Null
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null
}

function Set-ZeroTextCell($ref) {
    $ws.Range($ref).Value = "'0"
    $zeroTextSrc.Copy() | Out-Null
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null
}

function Set-NaTextCell($ref) {
    $ws.Range($ref).Value = "'***.*"
    $naTextSrc.Copy() | Out-Null
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null
}

# Row 16 - Murder
Set-ZeroTextCell "C16"
Set-NumCell "D16" 2
Set-PctCell "E16" -100
Set-NumCell "F16" 5
Set-NumCell "G16" 6
Set-PctCell "H16" -16.666666666666
Set-NumCell "J16" 6
Set-PctCell "K16" -16.666666666666
Set-PctCell "L16" 0
Set-PctCell "M16" 0
Set-PctCell "N16" -80

# Row 17 - Rape
Set-ZeroTextCell "C17"
Set-NumCell "D17" 1
Set-PctCell "E17" -100
Set-NumCell "F17" 4
Set-NumCell "G17" 8
Set-PctCell "H17" -50
Set-NumCell "J17" 8
Set-PctCell "K17" -50
Set-PctCell "M17" -33.333333333333
Set-PctCell "N17" -63.636363636363

# Row 18 - Robbery
Set-ZeroTextCell "C18"
Set-ZeroTextCell "D18"
Set-NaTextCell "E18"
Set-NumCell "F18" 8
Set-NumCell "G18" 4
Set-PctCell "H18" 100
Set-PctCell "L18" -12.5
Set-PctCell "M18" 0
Set-PctCell "N18" -80

# Row 19 - Fel. Assault
Set-NumCell "C19" 3
Set-NumCell "D19" 6
Set-NumCell "I19" 10
Set-NumCell "J19" 15
Set-PctCell "K19" -33.333333333333
Set-PctCell "L19" -33.333333333333
Set-PctCell "M19" -37.5
Set-PctCell "N19" -54.545454545454

# Row 20 - Burglary
Set-NumCell "D20" 1
Set-PctCell "E20" -100
Set-NumCell "F20" 4
Set-NumCell "G20" 5
Set-PctCell "H20" -20
Set-NumCell "J20" 5
Set-PctCell "K20" -20
Set-PctCell "L20" 0
Set-PctCell "M20" -33.333333333333
Set-PctCell "N20" -89.189189189189

# Row 21 - Gr. Larceny (bold/total styling, unaffected by style swap)
$ws.Range("C21").Value = 3
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = -70
$ws.Range("F21").Value = 33
$ws.Range("G21").Value = 39
$ws.Range("H21").Value = -15.384615384615
$ws.Range("I21").Value = 30
$ws.Range("J21").Value = 38
$ws.Range("K21").Value = -21.052631578947
$ws.Range("L21").Value = -21.052631578947
$ws.Range("M21").Value = -25
$ws.Range("N21").Value = -76.923076923076

# Row 22 - G.L.A.
Set-NumCell "C22" 1
Set-NumCell "F22" 1
Set-NumCell "I22" 1
Set-PctCell "M22" 0

# Row 23 - TOTAL
Set-ZeroTextCell "C23"
Set-NumCell "D23" 2
Set-PctCell "E23" -100
Set-NumCell "F23" 6
Set-NumCell "G23" 8
Set-PctCell "H23" -25
Set-NumCell "J23" 8
Set-PctCell "K23" -25
Set-PctCell "L23" -14.285714285714

# Row 24 - Transit
Set-NumCell "C24" 8
Set-NumCell "D24" 12
Set-PctCell "E24" -33.333333333333
Set-NumCell "F24" 40
Set-NumCell "G24" 29
Set-PctCell "H24" 37.931034482758
Set-NumCell "I24" 36
Set-NumCell "J24" 29
Set-PctCell "K24" 24.137931034482
Set-PctCell "L24" -12.195121951219
Set-PctCell "M24" 2.857142857142

# Row 25 - Housing
Set-NumCell "C25" 5
Set-NumCell "D25" 5
Set-PctCell "E25" 0
Set-NumCell "G25" 13
Set-PctCell "H25" 38.461538461538
Set-NumCell "I25" 15
Set-NumCell "J25" 13
Set-PctCell "K25" 15.384615384615
Set-PctCell "L25" 7.142857142857

# Row 26 - Petit Larceny
Set-NumCell "C26" 3
Set-NumCell "D26" 2
Set-PctCell "E26" 50
Set-NumCell "G26" 8
Set-PctCell "H26" 112.5
Set-NumCell "I26" 16
Set-NumCell "J26" 8
Set-PctCell "K26" 100
Set-PctCell "L26" -11.111111111111
Set-PctCell "M26" -30.434782608695

# Row 27 - Retail Theft
Set-ZeroTextCell "D27"
Set-NaTextCell "E27"

# Row 28 - Misd. Assault
Set-NumCell "D28" 1
Set-PctCell "E28" -100
Set-NumCell "G28" 1
Set-PctCell "H28" 0
Set-NumCell "J28" 1
Set-PctCell "K28" 0

# Row 31 - Shooting Inc.
Set-NumCell "F31" 1
Set-ZeroTextCell "G31"
Set-NaTextCell "H31"
Set-NumCell "I31" 1

# Row 33 - Traffic Fatalities
Set-ZeroTextCell "C33"
